# Applies the odds/match-data update for 2024-11-20 FlashScore workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 3.1
$ws.Range("I2").Value = 3.5
$ws.Range("J2").Value = 3.1
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5
$ws.Range("O2").Value = 1.57
$ws.Range("P2").Value = 2.25
$ws.Range("Q2").Value = 2.7
$ws.Range("R2").Value = 1.44
$ws.Range("S2").Value = 1.62
$ws.Range("T2").Value = 2.2
$ws.Range("U2").Value = 2.25
$ws.Range("V2").Value = 1.57
$ws.Range("Y2").Value = 10
$ws.Range("AA2").Value = 23
$ws.Range("AC2").Value = 6
$ws.Range("AF2").Value = 81
$ws.Range("AG2").Value = 7.5
$ws.Range("AI2").Value = 13
$ws.Range("AP2").Value = 29
$ws.Range("AR2").Value = 81
$ws.Range("AS2").Value = 351
$ws.Range("AT2").Value = 2.2
$ws.Range("AU2").Value = 9.5
$ws.Range("BA2").Value = 126
$ws.Range("BB2").Value = 401

# Row 4
$ws.Range("G4").Value = 1.7
$ws.Range("H4").Value = 3.75
$ws.Range("I4").Value = 5
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 2.15
$ws.Range("R4").Value = 1.67
$ws.Range("AD4").Value = 7
$ws.Range("AG4").Value = 11
$ws.Range("AH4").Value = 23
$ws.Range("AO4").Value = 9
$ws.Range("AS4").Value = 201
$ws.Range("AW4").Value = 6.5
$ws.Range("AZ4").Value = 101

# Row 7
$ws.Range("G7").Value = 1.5
$ws.Range("H7").Value = 4.1
$ws.Range("J7").Value = 2.1
$ws.Range("M7").Value = 1.07
$ws.Range("N7").Value = 9
$ws.Range("O7").Value = 1.36
$ws.Range("P7").Value = 3.2
$ws.Range("Q7").Value = 2.1
$ws.Range("R7").Value = 1.7
$ws.Range("Y7").Value = 9
$ws.Range("Z7").Value = 10
$ws.Range("AA7").Value = 15
$ws.Range("AC7").Value = 8.5
$ws.Range("AV7").Value = 81

# Row 8
$ws.Range("Q8").Value = 2.07
$ws.Range("R8").Value = 1.83

# Row 11
$ws.Range("G11").Value = 3.3
$ws.Range("H11").Value = 3.5
$ws.Range("I11").Value = 2.1
$ws.Range("J11").Value = 3.75
$ws.Range("L11").Value = 2.75
$ws.Range("M11").Value = 1.05
$ws.Range("N11").Value = 11
$ws.Range("O11").Value = 1.25
$ws.Range("P11").Value = 4
$ws.Range("U11").Value = 1.62
$ws.Range("V11").Value = 2.2
$ws.Range("Z11").Value = 34
$ws.Range("AA11").Value = 23
$ws.Range("AC11").Value = 12
$ws.Range("AD11").Value = 6.5
$ws.Range("AG11").Value = 9
$ws.Range("AH11").Value = 11
$ws.Range("AI11").Value = 9
$ws.Range("AJ11").Value = 19
$ws.Range("AO11").Value = 17
$ws.Range("AP11").Value = 23
$ws.Range("AR11").Value = 67
$ws.Range("AW11").Value = 4.33
$ws.Range("AX11").Value = 11
$ws.Range("AZ11").Value = 41

# Row 13
$ws.Range("Q13").Value = 1.93
$ws.Range("R13").Value = 1.97

# Row 14
$ws.Range("A14").Value = "IZ3qIEYa"
$ws.Range("C14").Value = "20:30"
$ws.Range("D14").Value = "COLOMBIA - PRIMERA A"
$ws.Range("E14").Value = "Millonarios"
$ws.Range("F14").Value = "Dep. Pasto"
$ws.Range("G14").Value = 1.85
$ws.Range("H14").Value = 3.1
$ws.Range("I14").Value = 4.75
$ws.Range("J14").Value = 2.63
$ws.Range("K14").Value = 1.91
$ws.Range("L14").Value = 5.5
$ws.Range("M14").Value = 1.11
$ws.Range("N14").Value = 6.5
$ws.Range("O14").Value = 1.5
$ws.Range("P14").Value = 2.5
$ws.Range("Q14").Value = 2.6
$ws.Range("R14").Value = 1.48
$ws.Range("S14").Value = 1.57
$ws.Range("T14").Value = 2.25
$ws.Range("U14").Value = 2.25
$ws.Range("V14").Value = 1.57
$ws.Range("W14").Value = 5
$ws.Range("X14").Value = 7.5
$ws.Range("Y14").Value = 9.5
$ws.Range("Z14").Value = 15
$ws.Range("AA14").Value = 19
$ws.Range("AB14").Value = 41
$ws.Range("AC14").Value = 6
$ws.Range("AD14").Value = 6.5
$ws.Range("AE14").Value = 21
$ws.Range("AF14").Value = 81
$ws.Range("AG14").Value = 9.5
$ws.Range("AH14").Value = 21
$ws.Range("AI14").Value = 17
$ws.Range("AJ14").Value = 51
$ws.Range("AK14").Value = 41
$ws.Range("AL14").Value = 51
$ws.Range("AM14").Value = 201
$ws.Range("AN14").Value = 3.6
$ws.Range("AO14").Value = 11
$ws.Range("AP14").Value = 29
$ws.Range("AQ14").Value = 41
$ws.Range("AR14").Value = 67
$ws.Range("AS14").Value = 301
$ws.Range("AT14").Value = 2.25
$ws.Range("AU14").Value = 10
$ws.Range("AV14").Value = 81
$ws.Range("AW14").Value = 6
$ws.Range("AX14").Value = 29
$ws.Range("AY14").Value = 41
$ws.Range("AZ14").Value = 101
$ws.Range("BA14").Value = 151
$ws.Range("BB14").Value = 501
$ws.Range("BC14").Value = 126
$ws.Range("BD14").Value = 126

# Row 16
$ws.Range("G16").Value = 35
$ws.Range("H16").Value = 8.25
$ws.Range("J16").Value = 24
$ws.Range("K16").Value = 3.35
$ws.Range("L16").Value = 1.3
$ws.Range("P16").Value = 6.2
$ws.Range("R16").Value = 3.25
$ws.Range("U16").Value = 2.77
$ws.Range("V16").Value = 1.39
$ws.Range("AB16").Value = 800
$ws.Range("AC16").Value = 20
$ws.Range("AF16").Value = 400
$ws.Range("AG16").Value = 10.5
$ws.Range("AL16").Value = 60
$ws.Range("AN16").Value = 32
$ws.Range("AO16").Value = 400
$ws.Range("AU16").Value = 14.5
$ws.Range("AX16").Value = 3.85
$ws.Range("AZ16").Value = 7.1
